$d = $word.ActiveDocument

$pairs = @(
    @("48×64=3072", "55×59=3245"),
    @("26×17=442",  "52×47=2444"),
    @("84×55=4620", "33×53=1749"),
    @("47×92=4324", "75×90=6750"),
    @("16×79=1264", "37×50=1850"),
    @("98×98=9604", "94×49=4606"),
    @("48×30=1440", "22×98=2156"),
    @("31×22=682",  "28×68=1904"),
    @("86×41=3526", "40×38=1520"),
    @("23×86=1978", "96×21=2016"),
    @("14×11=154",  "19×22=418"),
    @("72×40=2880", "94×34=3196"),
    @("38×66=2508", "52×48=2496"),
    @("70×17=1190", "82×95=7790"),
    @("81×44=3564", "68×28=1904"),
    @("94×84=7896", "13×25=325"),
    @("76×47=3572", "72×41=2952"),
    @("94×21=1974", "15×31=465"),
    @("11×55=605",  "90×45=4050"),
    @("75×30=2250", "94×86=8084"),
    @("23×77=1771", "32×98=3136"),
    @("42×67=2814", "55×26=1430"),
    @("90×33=2970", "53×73=3869"),
    @("34×91=3094", "36×44=1584"),
    @("96×28=2688", "25×60=1500")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
